# Daily attendance processing - 2025-10-29 21:43:06
# Normalize the "Recorded By" column (G) so that the current user
# (dnasr281@gmail.com) is always listed first in the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetUser = "dnasr281@gmail.com"

# Find the last used row in the sheet
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value.ToString().Contains(",") -and $value.ToString().Contains($targetUser)) {
        $parts = $value.ToString().Split(",")
        $trimmed = @()
        foreach ($p in $parts) { $trimmed += $p.Trim() }

        if ($trimmed.Count -eq 2 -and $trimmed[0] -ne $targetUser -and $trimmed[1] -eq $targetUser) {
            $newValue = $targetUser + ", " + $trimmed[0]
            $cell.Value = $newValue
        }
    }
}
